# Regenerate merged AHB files
# ---------------------------------------------------------------------------
# 1) Rename the header-row labels from the "_old"/"_new" suffix convention to
#    the version-specific "_FV2304"/"_FV2310" convention (the "diff" column,
#    K1, keeps its name).
# 2) Turn the data range A1:U54 into a real Excel Table ("Table1") so the new
#    column headers are picked up as the table's column names.
# 3) Freeze the header row (row 1) so it stays visible while scrolling.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells ------------------------------------------------
$headerNames = @{
    "A1" = "Segmentname_FV2304"
    "B1" = "Segmentgruppe_FV2304"
    "C1" = "Segment_FV2304"
    "D1" = "Datenelement_FV2304"
    "E1" = "Segment ID_FV2304"
    "F1" = "Code_FV2304"
    "G1" = "Qualifier_FV2304"
    "H1" = "Beschreibung_FV2304"
    "I1" = "Bedingungsausdruck_FV2304"
    "J1" = "Bedingung_FV2304"
    "L1" = "Segmentname_FV2310"
    "M1" = "Segmentgruppe_FV2310"
    "N1" = "Segment_FV2310"
    "O1" = "Datenelement_FV2310"
    "P1" = "Segment ID_FV2310"
    "Q1" = "Code_FV2310"
    "R1" = "Qualifier_FV2310"
    "S1" = "Beschreibung_FV2310"
    "T1" = "Bedingungsausdruck_FV2310"
    "U1" = "Bedingung_FV2310"
}

foreach ($addr in $headerNames.Keys) {
    $ws.Range($addr).Value = $headerNames[$addr]
}

# --- 2) Create the table over A1:U54 ---------------------------------------
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U54"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3) Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
